$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 46.595173
$ws.Cells.Item(2, 8).Value = 139.785519
$ws.Cells.Item(2, 9).Value = 0.7981698877785356
$ws.Cells.Item(2, 10).Value = 0.7981698877785355
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3547066666666667
$ws.Cells.Item(2, 14).Value = 1.06412
$ws.Cells.Item(2, 15).Value = 0.006934574868045491
$ws.Cells.Item(2, 16).Value = 0.00693457486804549
$ws.Cells.Item(2, 17).Value = 16.52761849758667
$ws.Cells.Item(2, 18).Value = 148.74856647828
$ws.Cells.Item(2, 19).Value = 0.005534968844219723
$ws.Cells.Item(2, 20).Value = 0.005534968844219721

# Row 3
$ws.Cells.Item(3, 7).Value = 46.595173
$ws.Cells.Item(3, 8).Value = 139.785519
$ws.Cells.Item(3, 9).Value = 0.7981698877785356
$ws.Cells.Item(3, 10).Value = 0.7981698877785355
$ws.Cells.Item(3, 15).Value = 0.01016563762403213
$ws.Cells.Item(3, 16).Value = 0.01016563762403213
$ws.Cells.Item(3, 17).Value = 24.228418271021
$ws.Cells.Item(3, 18).Value = 218.055764439189
$ws.Cells.Item(3, 19).Value = 0.008113905841570986
$ws.Cells.Item(3, 20).Value = 0.008113905841570982

# Row 4
$ws.Cells.Item(4, 7).Value = 46.595173
$ws.Cells.Item(4, 8).Value = 139.785519
$ws.Cells.Item(4, 9).Value = 0.7981698877785356
$ws.Cells.Item(4, 10).Value = 0.7981698877785355
$ws.Cells.Item(4, 13).Value = 26.50170333333334
$ws.Cells.Item(4, 14).Value = 79.50511
$ws.Cells.Item(4, 15).Value = 0.5181127482682332
$ws.Cells.Item(4, 16).Value = 0.5181127482682331
$ws.Cells.Item(4, 17).Value = 1234.851451611343
$ws.Cells.Item(4, 18).Value = 11113.66306450209
$ws.Cells.Item(4, 19).Value = 0.4135419941418843
$ws.Cells.Item(4, 20).Value = 0.4135419941418842

# Row 5
$ws.Cells.Item(5, 7).Value = 46.595173
$ws.Cells.Item(5, 8).Value = 139.785519
$ws.Cells.Item(5, 9).Value = 0.7981698877785356
$ws.Cells.Item(5, 10).Value = 0.7981698877785355
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.09443866666666667
$ws.Cells.Item(5, 14).Value = 0.283316
$ws.Cells.Item(5, 15).Value = 0.001846291784117559
$ws.Cells.Item(5, 16).Value = 0.001846291784117558
$ws.Cells.Item(5, 17).Value = 4.400386011222666
$ws.Cells.Item(5, 18).Value = 39.603474101004
$ws.Cells.Item(5, 19).Value = 0.001473654506135544
$ws.Cells.Item(5, 20).Value = 0.001473654506135544

# Row 6
$ws.Cells.Item(6, 7).Value = 46.595173
$ws.Cells.Item(6, 8).Value = 139.785519
$ws.Cells.Item(6, 9).Value = 0.7981698877785356
$ws.Cells.Item(6, 10).Value = 0.7981698877785355
$ws.Cells.Item(6, 13).Value = 23.67963033333333
$ws.Cells.Item(6, 14).Value = 71.03889099999999
$ws.Cells.Item(6, 15).Value = 0.4629407474555717
$ws.Cells.Item(6, 16).Value = 0.4629407474555717
$ws.Cells.Item(6, 17).Value = 1103.356471957714
$ws.Cells.Item(6, 18).Value = 9930.208247619428
$ws.Cells.Item(6, 19).Value = 0.3695053644447251
$ws.Cells.Item(6, 20).Value = 0.369505364444725

# Row 7
$ws.Cells.Item(7, 9).Value = 0.02837882113957134
$ws.Cells.Item(7, 10).Value = 0.02837882113957133
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3547066666666667
$ws.Cells.Item(7, 14).Value = 1.06412
$ws.Cells.Item(7, 15).Value = 0.006934574868045491
$ws.Cells.Item(7, 16).Value = 0.00693457486804549
$ws.Cells.Item(7, 17).Value = 0.5876372140666667
$ws.Cells.Item(7, 18).Value = 5.2887349266
$ws.Cells.Item(7, 19).Value = 0.0001967950598592295
$ws.Cells.Item(7, 20).Value = 0.0001967950598592294

# Row 8
$ws.Cells.Item(8, 9).Value = 0.02837882113957134
$ws.Cells.Item(8, 10).Value = 0.02837882113957133
$ws.Cells.Item(8, 15).Value = 0.01016563762403213
$ws.Cells.Item(8, 16).Value = 0.01016563762403213
$ws.Cells.Item(8, 19).Value = 0.0002884888119021048
$ws.Cells.Item(8, 20).Value = 0.0002884888119021047

# Row 9
$ws.Cells.Item(9, 9).Value = 0.02837882113957134
$ws.Cells.Item(9, 10).Value = 0.02837882113957133
$ws.Cells.Item(9, 13).Value = 26.50170333333334
$ws.Cells.Item(9, 14).Value = 79.50511
$ws.Cells.Item(9, 15).Value = 0.5181127482682332
$ws.Cells.Item(9, 16).Value = 0.5181127482682331
$ws.Cells.Item(9, 17).Value = 43.90497438678334
$ws.Cells.Item(9, 18).Value = 395.14476948105
$ws.Cells.Item(9, 19).Value = 0.01470342901323594
$ws.Cells.Item(9, 20).Value = 0.01470342901323593

# Row 10
$ws.Cells.Item(10, 9).Value = 0.02837882113957134
$ws.Cells.Item(10, 10).Value = 0.02837882113957133
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.09443866666666667
$ws.Cells.Item(10, 14).Value = 0.283316
$ws.Cells.Item(10, 15).Value = 0.001846291784117559
$ws.Cells.Item(10, 16).Value = 0.001846291784117558
$ws.Cells.Item(10, 17).Value = 0.1564551224866667
$ws.Cells.Item(10, 18).Value = 1.40809610238
$ws.Cells.Item(10, 19).Value = 0.00005239558431293225
$ws.Cells.Item(10, 20).Value = 0.00005239558431293223

# Row 11
$ws.Cells.Item(11, 9).Value = 0.02837882113957134
$ws.Cells.Item(11, 10).Value = 0.02837882113957133
$ws.Cells.Item(11, 13).Value = 23.67963033333333
$ws.Cells.Item(11, 14).Value = 71.03889099999999
$ws.Cells.Item(11, 15).Value = 0.4629407474555717
$ws.Cells.Item(11, 16).Value = 0.4629407474555717
$ws.Cells.Item(11, 17).Value = 39.22968837877833
$ws.Cells.Item(11, 18).Value = 353.067195409005
$ws.Cells.Item(11, 19).Value = 0.01313771267026113
$ws.Cells.Item(11, 20).Value = 0.01313771267026113

# Row 12
$ws.Cells.Item(12, 7).Value = 5.966798333333333
$ws.Cells.Item(12, 8).Value = 17.900395
$ws.Cells.Item(12, 9).Value = 0.102210560654294
$ws.Cells.Item(12, 10).Value = 0.1022105606542939
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.3547066666666667
$ws.Cells.Item(12, 14).Value = 1.06412
$ws.Cells.Item(12, 15).Value = 0.006934574868045491
$ws.Cells.Item(12, 16).Value = 0.00693457486804549
$ws.Cells.Item(12, 17).Value = 2.116463147488889
$ws.Cells.Item(12, 18).Value = 19.0481683274
$ws.Cells.Item(12, 19).Value = 0.0007087867851621062
$ws.Cells.Item(12, 20).Value = 0.000708786785162106

# Row 13
$ws.Cells.Item(13, 7).Value = 5.966798333333333
$ws.Cells.Item(13, 8).Value = 17.900395
$ws.Cells.Item(13, 9).Value = 0.102210560654294
$ws.Cells.Item(13, 10).Value = 0.1022105606542939
$ws.Cells.Item(13, 15).Value = 0.01016563762403213
$ws.Cells.Item(13, 16).Value = 0.01016563762403213
$ws.Cells.Item(13, 17).Value = 3.102597896971667
$ws.Cells.Item(13, 18).Value = 27.923381072745
$ws.Cells.Item(13, 19).Value = 0.001039035520960709
$ws.Cells.Item(13, 20).Value = 0.001039035520960709

# Row 14
$ws.Cells.Item(14, 7).Value = 5.966798333333333
$ws.Cells.Item(14, 8).Value = 17.900395
$ws.Cells.Item(14, 9).Value = 0.102210560654294
$ws.Cells.Item(14, 10).Value = 0.1022105606542939
$ws.Cells.Item(14, 13).Value = 26.50170333333334
$ws.Cells.Item(14, 14).Value = 79.50511
$ws.Cells.Item(14, 15).Value = 0.5181127482682332
$ws.Cells.Item(14, 16).Value = 0.5181127482682331
$ws.Cells.Item(14, 17).Value = 158.1303192798278
$ws.Cells.Item(14, 18).Value = 1423.17287351845
$ws.Cells.Item(14, 19).Value = 0.05295659448263319
$ws.Cells.Item(14, 20).Value = 0.05295659448263317

# Row 15
$ws.Cells.Item(15, 7).Value = 5.966798333333333
$ws.Cells.Item(15, 8).Value = 17.900395
$ws.Cells.Item(15, 9).Value = 0.102210560654294
$ws.Cells.Item(15, 10).Value = 0.1022105606542939
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.09443866666666667
$ws.Cells.Item(15, 14).Value = 0.283316
$ws.Cells.Item(15, 15).Value = 0.001846291784117559
$ws.Cells.Item(15, 16).Value = 0.001846291784117558
$ws.Cells.Item(15, 17).Value = 0.5634964788688889
$ws.Cells.Item(15, 18).Value = 5.07146830982
$ws.Cells.Item(15, 19).Value = 0.0001887105183860723
$ws.Cells.Item(15, 20).Value = 0.0001887105183860723

# Row 16
$ws.Cells.Item(16, 7).Value = 5.966798333333333
$ws.Cells.Item(16, 8).Value = 17.900395
$ws.Cells.Item(16, 9).Value = 0.102210560654294
$ws.Cells.Item(16, 10).Value = 0.1022105606542939
$ws.Cells.Item(16, 13).Value = 23.67963033333333
$ws.Cells.Item(16, 14).Value = 71.03889099999999
$ws.Cells.Item(16, 15).Value = 0.4629407474555717
$ws.Cells.Item(16, 16).Value = 0.4629407474555717
$ws.Cells.Item(16, 17).Value = 141.2915788068828
$ws.Cells.Item(16, 18).Value = 1271.624209261945
$ws.Cells.Item(16, 19).Value = 0.0473174333471519
$ws.Cells.Item(16, 20).Value = 0.04731743334715188

# Row 17
$ws.Cells.Item(17, 7).Value = 0.4842143333333333
$ws.Cells.Item(17, 8).Value = 1.452643
$ws.Cells.Item(17, 9).Value = 0.008294535146321381
$ws.Cells.Item(17, 10).Value = 0.008294535146321381
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.3547066666666667
$ws.Cells.Item(17, 14).Value = 1.06412
$ws.Cells.Item(17, 15).Value = 0.006934574868045491
$ws.Cells.Item(17, 16).Value = 0.00693457486804549
$ws.Cells.Item(17, 17).Value = 0.1717540521288889
$ws.Cells.Item(17, 18).Value = 1.54578646916
$ws.Cells.Item(17, 19).Value = 0.00005751907496780028
$ws.Cells.Item(17, 20).Value = 0.00005751907496780027

# Row 18
$ws.Cells.Item(18, 7).Value = 0.4842143333333333
$ws.Cells.Item(18, 8).Value = 1.452643
$ws.Cells.Item(18, 9).Value = 0.008294535146321381
$ws.Cells.Item(18, 10).Value = 0.008294535146321381
$ws.Cells.Item(18, 15).Value = 0.01016563762403213
$ws.Cells.Item(18, 16).Value = 0.01016563762403213
$ws.Cells.Item(18, 17).Value = 0.2517803164036667
$ws.Cells.Item(18, 18).Value = 2.266022847633
$ws.Cells.Item(18, 19).Value = 0.00008431923855730149
$ws.Cells.Item(18, 20).Value = 0.00008431923855730148

# Row 19
$ws.Cells.Item(19, 7).Value = 0.4842143333333333
$ws.Cells.Item(19, 8).Value = 1.452643
$ws.Cells.Item(19, 9).Value = 0.008294535146321381
$ws.Cells.Item(19, 10).Value = 0.008294535146321381
$ws.Cells.Item(19, 13).Value = 26.50170333333334
$ws.Cells.Item(19, 14).Value = 79.50511
$ws.Cells.Item(19, 15).Value = 0.5181127482682332
$ws.Cells.Item(19, 16).Value = 0.5181127482682331
$ws.Cells.Item(19, 17).Value = 12.83250461174778
$ws.Cells.Item(19, 18).Value = 115.49254150573
$ws.Cells.Item(19, 19).Value = 0.004297504400268022
$ws.Cells.Item(19, 20).Value = 0.004297504400268022

# Row 20
$ws.Cells.Item(20, 7).Value = 0.4842143333333333
$ws.Cells.Item(20, 8).Value = 1.452643
$ws.Cells.Item(20, 9).Value = 0.008294535146321381
$ws.Cells.Item(20, 10).Value = 0.008294535146321381
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 0.09443866666666667
$ws.Cells.Item(20, 14).Value = 0.283316
$ws.Cells.Item(20, 15).Value = 0.001846291784117559
$ws.Cells.Item(20, 16).Value = 0.001846291784117558
$ws.Cells.Item(20, 17).Value = 0.04572855602088889
$ws.Cells.Item(20, 18).Value = 0.411557004188
$ws.Cells.Item(20, 19).Value = 0.0000153141320937275
$ws.Cells.Item(20, 20).Value = 0.0000153141320937275

# Row 21
$ws.Cells.Item(21, 7).Value = 0.4842143333333333
$ws.Cells.Item(21, 8).Value = 1.452643
$ws.Cells.Item(21, 9).Value = 0.008294535146321381
$ws.Cells.Item(21, 10).Value = 0.008294535146321381
$ws.Cells.Item(21, 13).Value = 23.67963033333333
$ws.Cells.Item(21, 14).Value = 71.03889099999999
$ws.Cells.Item(21, 15).Value = 0.4629407474555717
$ws.Cells.Item(21, 16).Value = 0.4629407474555717
$ws.Cells.Item(21, 17).Value = 11.46601641543478
$ws.Cells.Item(21, 18).Value = 103.194147738913
$ws.Cells.Item(21, 19).Value = 0.00383987830043453
$ws.Cells.Item(21, 20).Value = 0.00383987830043453

# Row 22
$ws.Cells.Item(22, 7).Value = 3.674642333333333
$ws.Cells.Item(22, 8).Value = 11.023927
$ws.Cells.Item(22, 9).Value = 0.0629461952812778
$ws.Cells.Item(22, 10).Value = 0.0629461952812778
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 0.3547066666666667
$ws.Cells.Item(22, 14).Value = 1.06412
$ws.Cells.Item(22, 15).Value = 0.006934574868045491
$ws.Cells.Item(22, 16).Value = 0.00693457486804549
$ws.Cells.Item(22, 17).Value = 1.303420133248889
$ws.Cells.Item(22, 18).Value = 11.73078119924
$ws.Cells.Item(22, 19).Value = 0.0004365051038366327
$ws.Cells.Item(22, 20).Value = 0.0004365051038366327

# Row 23
$ws.Cells.Item(23, 7).Value = 3.674642333333333
$ws.Cells.Item(23, 8).Value = 11.023927
$ws.Cells.Item(23, 9).Value = 0.0629461952812778
$ws.Cells.Item(23, 10).Value = 0.0629461952812778
$ws.Cells.Item(23, 15).Value = 0.01016563762403213
$ws.Cells.Item(23, 16).Value = 0.01016563762403213
$ws.Cells.Item(23, 17).Value = 1.910729496559667
$ws.Cells.Item(23, 18).Value = 17.196565469037
$ws.Cells.Item(23, 19).Value = 0.0006398882110410315
$ws.Cells.Item(23, 20).Value = 0.0006398882110410314

# Row 24
$ws.Cells.Item(24, 7).Value = 3.674642333333333
$ws.Cells.Item(24, 8).Value = 11.023927
$ws.Cells.Item(24, 9).Value = 0.0629461952812778
$ws.Cells.Item(24, 10).Value = 0.0629461952812778
$ws.Cells.Item(24, 13).Value = 26.50170333333334
$ws.Cells.Item(24, 14).Value = 79.50511
$ws.Cells.Item(24, 15).Value = 0.5181127482682332
$ws.Cells.Item(24, 16).Value = 0.5181127482682331
$ws.Cells.Item(24, 17).Value = 97.38428097410778
$ws.Cells.Item(24, 18).Value = 876.4585287669701
$ws.Cells.Item(24, 19).Value = 0.03261322623021173
$ws.Cells.Item(24, 20).Value = 0.03261322623021173

# Row 25
$ws.Cells.Item(25, 7).Value = 3.674642333333333
$ws.Cells.Item(25, 8).Value = 11.023927
$ws.Cells.Item(25, 9).Value = 0.0629461952812778
$ws.Cells.Item(25, 10).Value = 0.0629461952812778
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 0.09443866666666667
$ws.Cells.Item(25, 14).Value = 0.283316
$ws.Cells.Item(25, 15).Value = 0.001846291784117559
$ws.Cells.Item(25, 16).Value = 0.001846291784117558
$ws.Cells.Item(25, 17).Value = 0.3470283224368889
$ws.Cells.Item(25, 18).Value = 3.123254901932
$ws.Cells.Item(25, 19).Value = 0.0001162170431892826
$ws.Cells.Item(25, 20).Value = 0.0001162170431892826

# Row 26
$ws.Cells.Item(26, 7).Value = 3.674642333333333
$ws.Cells.Item(26, 8).Value = 11.023927
$ws.Cells.Item(26, 9).Value = 0.0629461952812778
$ws.Cells.Item(26, 10).Value = 0.0629461952812778
$ws.Cells.Item(26, 13).Value = 23.67963033333333
$ws.Cells.Item(26, 14).Value = 71.03889099999999
$ws.Cells.Item(26, 15).Value = 0.4629407474555717
$ws.Cells.Item(26, 16).Value = 0.4629407474555717
$ws.Cells.Item(26, 17).Value = 87.01417206055078
$ws.Cells.Item(26, 18).Value = 783.127548544957
$ws.Cells.Item(26, 19).Value = 0.02914035869299912
$ws.Cells.Item(26, 20).Value = 0.02914035869299912
